$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 634.26666
$ws.Range("I5").Value = 411.9
$ws.Range("J5").Value = 1079
$ws.Range("K5").Value = 411.9
$ws.Range("L5").Value = 1079
$ws.Range("M5").Value = -296.9
$ws.Range("N5").Value = -1309

$ws.Range("H12").Value = 196.4
$ws.Range("I12").Value = 222.5
$ws.Range("J12").Value = 92
$ws.Range("K12").Value = 222.5
$ws.Range("L12").Value = 92
$ws.Range("M12").Value = -52.5
$ws.Range("N12").Value = -432

$ws.Range("H28").Value = 245.42857
$ws.Range("I28").Value = 143.6
$ws.Range("J28").Value = 500
$ws.Range("K28").Value = 143.6
$ws.Range("L28").Value = 500
$ws.Range("M28").Value = 341.4
$ws.Range("N28").Value = -1470

$ws.Range("H69").Value = 134628.75
$ws.Range("I69").Value = 8500
$ws.Range("J69").Value = 176671.67
$ws.Range("K69").Value = 25500
$ws.Range("L69").Value = 530015.01
$ws.Range("M69").Value = -24626
$ws.Range("N69").Value = -531763.01

$ws.Range("H72").Value = 134628.75
$ws.Range("I72").Value = 8500
$ws.Range("J72").Value = 176671.67
$ws.Range("K72").Value = 76500
$ws.Range("L72").Value = 1590045.03
$ws.Range("M72").Value = -72132
$ws.Range("N72").Value = -1598781.03

$ws.Range("H76").Value = 1933.3334
$ws.Range("I76").Value = 1400
$ws.Range("K76").Value = 1400
$ws.Range("M76").Value = -1085

$ws.Range("H79").Value = 1933.3334
$ws.Range("I79").Value = 1400
$ws.Range("K79").Value = 1400
$ws.Range("M79").Value = -308

$ws.Range("H116").Value = 4400
$ws.Range("I116").Value = 4247
$ws.Range("J116").Value = 4706
$ws.Range("K116").Value = 4247
$ws.Range("L116").Value = 4706
$ws.Range("M116").Value = -805
$ws.Range("N116").Value = -11590

$ws.Range("H137").Value = 3129.3684
$ws.Range("I137").Value = 2879.8823
$ws.Range("K137").Value = 8639.6469
$ws.Range("M137").Value = -6089.6469

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 6218
$ws.Range("I6").Value = 545
$ws.Range("J6").Value = 10000
$ws.Range("K6").Value = 545
$ws.Range("L6").Value = 10000
$ws.Range("M6").Value = -372
$ws.Range("N6").Value = -10346

$ws.Range("H32").Value = 7216.9
$ws.Range("I32").Value = 7216.9
$ws.Range("K32").Value = 7216.9
$ws.Range("M32").Value = -6929.9

$ws.Range("H45").Value = 4302.8335
$ws.Range("I45").Value = 4580.6
$ws.Range("K45").Value = 4580.6
$ws.Range("M45").Value = -4203.6

$ws.Range("H74").Value = 30668.334
$ws.Range("I74").Value = 28802.2
$ws.Range("K74").Value = 28802.2
$ws.Range("M74").Value = -27928.2

$ws.Range("H77").Value = 30668.334
$ws.Range("I77").Value = 28802.2
$ws.Range("K77").Value = 144011
$ws.Range("M77").Value = -139643

$ws.Range("H101").Value = 100000
$ws.Range("J101").Value = 100000
$ws.Range("L101").Value = 100000
$ws.Range("N101").Value = -106490

$ws.Range("H102").Value = 2834.3333
$ws.Range("I102").Value = 2637.4546
$ws.Range("K102").Value = 2637.4546
$ws.Range("M102").Value = -1015.4546

$ws.Range("H132").Value = 2015.08
$ws.Range("I132").Value = 1470.4762
$ws.Range("K132").Value = 4411.4286
$ws.Range("M132").Value = -1881.4286

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 975.5
$ws.Range("I11").Value = 463.25
$ws.Range("K11").Value = 463.25
$ws.Range("M11").Value = -323.25

$ws.Range("H36").Value = 3282.6667
$ws.Range("I36").Value = 3282.6667
$ws.Range("K36").Value = 3282.6667
$ws.Range("M36").Value = -2748.6667

$ws.Range("H105").Value = 2851.7693
$ws.Range("I105").Value = 2851.7693
$ws.Range("K105").Value = 2851.7693
$ws.Range("M105").Value = -1104.7693

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 287999.72
$ws.Range("I3").Value = 501249.75
$ws.Range("J3").Value = 3666.3333
$ws.Range("K3").Value = 501249.75
$ws.Range("L3").Value = 3666.3333
$ws.Range("M3").Value = -501136.75
$ws.Range("N3").Value = -3892.3333

$ws.Range("H6").Value = 5375059
$ws.Range("I6").Value = 10750000
$ws.Range("J6").Value = 117.5
$ws.Range("K6").Value = 10750000
$ws.Range("L6").Value = 117.5
$ws.Range("M6").Value = -10749887
$ws.Range("N6").Value = -343.5

$ws.Range("H7").Value = 89.55556
$ws.Range("I7").Value = 100
$ws.Range("J7").Value = 68.666664
$ws.Range("K7").Value = 100
$ws.Range("L7").Value = 68.666664
$ws.Range("M7").Value = 13
$ws.Range("N7").Value = -294.666664

$ws.Range("H10").Value = 3332.5715
$ws.Range("I10").Value = 556.6667
$ws.Range("J10").Value = 19988
$ws.Range("K10").Value = 556.6667
$ws.Range("L10").Value = 19988
$ws.Range("M10").Value = -417.6667
$ws.Range("N10").Value = -20266

$ws.Range("H25").Value = 7836.8335
$ws.Range("I25").Value = 1255.25
$ws.Range("J25").Value = 21000
$ws.Range("K25").Value = 1255.25
$ws.Range("L25").Value = 21000
$ws.Range("M25").Value = -1081.25
$ws.Range("N25").Value = -21348

$ws.Range("H31").Value = 1286.6875
$ws.Range("I31").Value = 1014.46155
$ws.Range("J31").Value = 2466.3333
$ws.Range("K31").Value = 1014.46155
$ws.Range("L31").Value = 2466.3333
$ws.Range("M31").Value = -719.46155
$ws.Range("N31").Value = -3056.3333

$ws.Range("H34").Value = 1286.6875
$ws.Range("I34").Value = 1014.46155
$ws.Range("J34").Value = 2466.3333
$ws.Range("K34").Value = 1014.46155
$ws.Range("L34").Value = 2466.3333
$ws.Range("M34").Value = -812.46155
$ws.Range("N34").Value = -2870.3333

$ws.Range("H41").Value = 5039
$ws.Range("I41").Value = 5039
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 5039
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = -4611
$ws.Range("N41").ClearContents()

$ws.Range("H58").Value = 3349.5
$ws.Range("I58").Value = 3349.5
$ws.Range("K58").Value = 3349.5
$ws.Range("M58").Value = -3146.5

$ws.Range("H132").Value = 1282.591
$ws.Range("I132").Value = 910.95
$ws.Range("K132").Value = 2732.85
$ws.Range("M132").Value = -202.8500000000004

$ws.Range("H134").Value = 2961.2632
$ws.Range("I134").Value = 2570.2222
$ws.Range("K134").Value = 7710.6666
$ws.Range("M134").Value = -5175.6666

$ws.Range("H136").Value = 3349.5
$ws.Range("I136").Value = 3349.5
$ws.Range("K136").Value = 10048.5
$ws.Range("M136").Value = -7498.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 10016.5
$ws.Range("I4").Value = 10016.5
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 30049.5
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -29937.5
$ws.Range("N4").ClearContents()

$ws.Range("H7").Value = 50
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()

$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 350050000
$ws.Range("J2").Value = 99999
$ws.Range("L2").Value = 99999
$ws.Range("N2").Value = -100223

$ws.Range("H40").Value = 2899.8572
$ws.Range("I40").Value = 2967.3333
$ws.Range("J40").Value = 2495
$ws.Range("K40").Value = 2967.3333
$ws.Range("L40").Value = 2495
$ws.Range("M40").Value = -2831.3333
$ws.Range("N40").Value = -2767

$ws.Range("H132").Value = 3024
$ws.Range("I132").Value = 2513.4707
$ws.Range("K132").Value = 7540.4121
$ws.Range("M132").Value = -5010.4121

$ws.Range("H136").Value = 4666.3335
$ws.Range("I136").Value = 4666.3335
$ws.Range("K136").Value = 13999.0005
$ws.Range("M136").Value = -11449.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 264966.25
$ws.Range("I2").Value = 264966.25
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 264966.25
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -264854.25
$ws.Range("N2").ClearContents()

$ws.Range("H122").Value = 2724.7
$ws.Range("I122").Value = 2805.5293
$ws.Range("J122").Value = 2266.6667
$ws.Range("K122").Value = 8416.5879
$ws.Range("L122").Value = 6800.000100000001
$ws.Range("M122").Value = -5966.5879
$ws.Range("N122").Value = -11700.0001

$ws.Range("H132").Value = 1443.1428
$ws.Range("I132").Value = 836.9091
$ws.Range("K132").Value = 2510.7273
$ws.Range("M132").Value = 19.27269999999999

$ws.Range("H136").Value = 2132.5667
$ws.Range("I136").Value = 1943.7778
$ws.Range("K136").Value = 5831.3334
$ws.Range("M136").Value = -3281.3334
